$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temporarily format the Price/Volume data range as Text so that values
# like "29.937.15" or "0.9999" are stored as strings, not reinterpreted as numbers.
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

$ws.Range('D2').Value = '29.937.15'
$ws.Range('E2').Value = '  +1.56%  '
$ws.Range('D3').Value = '1.864.03'
$ws.Range('E3').Value = '  +1.26%  '
$ws.Range('D4').Value = '0.9999'
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').Value = '244.96'
$ws.Range('E5').Value = '  +0.53%  '
$ws.Range('D6').Value = '0.6608'
$ws.Range('E6').Value = '  +5.39%  '
$ws.Range('D7').Value = '1.001'
$ws.Range('D8').Value = '0.07584'
$ws.Range('E8').Value = '  +2.23%  '
$ws.Range('D9').Value = '0.3006'
$ws.Range('E9').Value = '  +1.79%  '
$ws.Range('D10').Value = '24.78'
$ws.Range('D11').Value = '0.07673'
$ws.Range('E11').Value = '  -0.11%  '
$ws.Range('D12').Value = '1.871.46'
$ws.Range('E12').Value = '  +1.29%  '
$ws.Range('D13').Value = '5.089'
$ws.Range('E13').Value = '  +1.15%  '
$ws.Range('D14').Value = '0.6956'
$ws.Range('E14').Value = '  +2.62%  '
$ws.Range('D15').Value = '84.11'
$ws.Range('E15').Value = '  +0.91%  '
$ws.Range('D16').Value = '0.000009688'
$ws.Range('E16').Value = '  +5.77%  '
$ws.Range('D17').Value = '6.163'
$ws.Range('E17').Value = '  +4.21%  '
$ws.Range('D18').Value = '29.959.37'
$ws.Range('E18').Value = '  +1.56%  '
$ws.Range('D19').Value = '2.122.09'
$ws.Range('E19').Value = '  +1.24%  '
$ws.Range('D20').Value = '237.06'
$ws.Range('E20').Value = '  -2.38%  '
$ws.Range('D21').Value = '12.74'
$ws.Range('E21').Value = '  +1.29%  '
$ws.Range('E22').Value = '  -0.07%  '
$ws.Range('D23').Value = '7.772'
$ws.Range('E23').Value = '  +4.67%  '
$ws.Range('E24').Value = '  -0.09%  '
$ws.Range('D25').Value = '0.1454'
$ws.Range('E25').Value = '  +2.95%  '
$ws.Range('D26').Value = '159.17'
$ws.Range('E26').Value = '  +0.36%  '
$ws.Range('D27').Value = '8.629'
$ws.Range('E27').Value = '  +1.01%  '
$ws.Range('D28').Value = '17.98'
$ws.Range('E28').Value = '  +1.04%  '
$ws.Range('D29').Value = '0.06086'
$ws.Range('E29').Value = '  -1.31%  '
$ws.Range('D30').Value = '1.499'
$ws.Range('E30').Value = '  +0.13%  '
$ws.Range('E31').Value = '  +5.32%  '
$ws.Range('D32').Value = '4.179'
$ws.Range('E32').Value = '  +1.25%  '
$ws.Range('D33').Value = '4.122'
$ws.Range('E33').Value = '  +0.54%  '
$ws.Range('D34').Value = '1.888'
$ws.Range('E34').Value = '  +0.65%  '
$ws.Range('D35').Value = '1.183'
$ws.Range('E35').Value = '  +3.36%  '
$ws.Range('D36').Value = '0.7400'
$ws.Range('E36').Value = '  +1.60%  '
$ws.Range('D38').Value = '2.819'
$ws.Range('E38').Value = '  -2.55%  '
$ws.Range('D39').Value = '0.01801'
$ws.Range('E39').Value = '  +1.90%  '
$ws.Range('D40').Value = '1.217.69'
$ws.Range('E40').Value = '  -0.76%  '
$ws.Range('D41').Value = '6.399'
$ws.Range('E41').Value = '  +1.14%  '
$ws.Range('D42').Value = '0.9160'
$ws.Range('E42').Value = '  +0.15%  '
$ws.Range('D43').Value = '2.032.82'
$ws.Range('E43').Value = '  +1.15%  '
$ws.Range('E44').Value = '  -0.15%  '
$ws.Range('D45').Value = '7.774'
$ws.Range('E45').Value = '  +16.35%  '
$ws.Range('D46').Value = '67.78'
$ws.Range('E46').Value = '  +3.06%  '
$ws.Range('D47').Value = '101.83'
$ws.Range('E47').Value = '  -0.05%  '
$ws.Range('E48').Value = '  -0.82%  '
$ws.Range('D49').Value = '0.4095'
$ws.Range('E49').Value = '  +0.73%  '
$ws.Range('D50').Value = '9.216'
$ws.Range('E50').Value = '  -0.16%  '
$ws.Range('D51').Value = '1.692'
$ws.Range('E51').Value = '  +5.43%  '

# Restore the original (default) cell style now that the text values are set.
$dataRange.Style = "Normal"
